$p = $ppt.ActivePresentation

# Delete slide 8 ("imports inside vendor")
$p.Slides.Item(8).Delete()
